# Adjusted graph positions in ppt
#
# Slide 4 ("CNN with 2 Convolutional layers ...") holds a title text box and
# three picture shapes (the compiled graphs). They are resized/repositioned
# so the pictures are bigger and better centred on the slide.
#
# NOTE on numeric literals: PowerPoint's Shape.Left/Top/Width/Height
# properties are exposed as single-precision (32-bit) floats measured in
# points, while the underlying OOXML stores EMUs (1 pt = 12700 EMU) as
# integers. A naive `emu / 12700.0` can therefore round-trip to an EMU
# value that is off by one once it has been truncated to float32. The
# literals below were chosen so that, after the float32 round-trip, they
# land exactly on the intended EMU values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Shape 1: "TextBox 1" - title text box
$title = $s.Shapes.Item(1)
$title.Left   = 28.54051244094488
$title.Top    = 12.815473070866142
$title.Width  = 902.9188538976379
$title.Height = 41.198464196850395

# Shape 2: "Picture 3" - first graph (top-left)
$pic1 = $s.Shapes.Item(2)
$pic1.Left   = 4.466339212598426
$pic1.Top    = 61.29736383464567
$pic1.Width  = 383.13870241732286
$pic1.Height = 255.42578203149606

# Shape 3: "Picture 5" - second graph (top-right)
$pic2 = $s.Shapes.Item(3)
$pic2.Left   = 572.3950196299212
$pic2.Top    = 61.29736383464567
$pic2.Width  = 383.13870241732286
$pic2.Height = 255.42578203149606

# Shape 4: "Picture 7" - third graph (bottom-centre)
$pic3 = $s.Shapes.Item(4)
$pic3.Left   = 288.4305114409449
$pic3.Top    = 270.000031
$pic3.Width  = 383.13870241732286
$pic3.Height = 255.42578203149606
